$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove deleted samples (RM 232 at original row 26, SC 92 at original row 28)
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# Apply updated / re-imputed cell values (post row-removal row numbers)
$ws.Range("E2").Value = -7.2
$ws.Range("F3").ClearContents()
$ws.Range("F4").Value = 17.97
$ws.Range("F5").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("E12").Value = -5.3
$ws.Range("E14").ClearContents()
$ws.Range("F15").Value = 16.2
$ws.Range("F18").Value = 18.35
$ws.Range("F19").ClearContents()
$ws.Range("E20").Value = -7.2
$ws.Range("E21").Value = -8.7
$ws.Range("F22").ClearContents()
$ws.Range("E23").ClearContents()
$ws.Range("F23").Value = 16.48
$ws.Range("E24").ClearContents()
$ws.Range("F25").Value = 16.6
$ws.Range("C26").Value = 10.8
$ws.Range("C27").ClearContents()
$ws.Range("F27").ClearContents()
$ws.Range("C30").Value = 11.4
$ws.Range("E31").Value = -8.1
$ws.Range("C32").ClearContents()
$ws.Range("E33").Value = -10.7

Write-Output "edit applied"
